$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("C3").Value = "Desenho Técnico - MCT-1A"
$ws.Range("E3").Value = "Desenho Técnico - MEC-1A"
$ws.Range("F3").Value = "CAD - MCT-2A"

# Row 4 updates
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "CAD - MEC-2A"
$ws.Range("F4").Value = "-"

# Row 6 updates
$ws.Range("D6").Value = "Desenho Técnico - ELT-1A"
